$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 4437.067
$ws.Cells.Item(40, 10).Value = 4572.2
$ws.Cells.Item(40, 12).Value = 4572.2
$ws.Cells.Item(40, 14).Value = -4922.2
$ws.Cells.Item(86, 8).Value = 44570.3
$ws.Cells.Item(86, 9).Value = 1534.3334
$ws.Cells.Item(86, 10).Value = 63014.285
$ws.Cells.Item(86, 11).Value = 1534.3334
$ws.Cells.Item(86, 12).Value = 63014.285
$ws.Cells.Item(86, 13).Value = -411.3334
$ws.Cells.Item(86, 14).Value = -65260.285
$ws.Cells.Item(89, 8).Value = 44570.3
$ws.Cells.Item(89, 9).Value = 1534.3334
$ws.Cells.Item(89, 10).Value = 63014.285
$ws.Cells.Item(89, 11).Value = 7671.666999999999
$ws.Cells.Item(89, 12).Value = 315071.425
$ws.Cells.Item(89, 13).Value = -2055.666999999999
$ws.Cells.Item(89, 14).Value = -326303.425
$ws.Cells.Item(94, 8).Value = 38561236
$ws.Cells.Item(94, 9).Value = 50006604
$ws.Cells.Item(94, 11).Value = 50006604
$ws.Cells.Item(94, 13).Value = -50006153
$ws.Cells.Item(138, 8).Value = 4500.033
$ws.Cells.Item(138, 9).Value = 1191.0769
$ws.Cells.Item(138, 10).Value = 5051.526
$ws.Cells.Item(138, 11).Value = 3573.2307
$ws.Cells.Item(138, 12).Value = 15154.578
$ws.Cells.Item(138, 13).Value = 1566.7693
$ws.Cells.Item(138, 14).Value = -25434.578

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 8727.223
$ws.Cells.Item(2, 9).Value = 11035.667
$ws.Cells.Item(2, 10).Value = 4110.3335
$ws.Cells.Item(2, 11).Value = 11035.667
$ws.Cells.Item(2, 12).Value = 4110.3335
$ws.Cells.Item(2, 13).Value = -10922.667
$ws.Cells.Item(2, 14).Value = -4336.3335
$ws.Cells.Item(32, 8).Value = 2360.24
$ws.Cells.Item(32, 9).Value = 2291.9167
$ws.Cells.Item(32, 11).Value = 2291.9167
$ws.Cells.Item(32, 13).Value = -2004.9167
$ws.Cells.Item(116, 8).Value = 8727.223
$ws.Cells.Item(116, 9).Value = 11035.667
$ws.Cells.Item(116, 10).Value = 4110.3335
$ws.Cells.Item(116, 11).Value = 11035.667
$ws.Cells.Item(116, 12).Value = 4110.3335
$ws.Cells.Item(116, 13).Value = -8741.666999999999
$ws.Cells.Item(116, 14).Value = -8698.333500000001
$ws.Cells.Item(132, 8).Value = 2486.6099
$ws.Cells.Item(132, 9).Value = 1651.8387
$ws.Cells.Item(132, 11).Value = 4955.5161
$ws.Cells.Item(132, 13).Value = -2425.5161

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 8727.223
$ws.Cells.Item(3, 9).Value = 11035.667
$ws.Cells.Item(3, 10).Value = 4110.3335
$ws.Cells.Item(3, 11).Value = 11035.667
$ws.Cells.Item(3, 12).Value = 4110.3335
$ws.Cells.Item(3, 13).Value = -10921.667
$ws.Cells.Item(3, 14).Value = -4338.3335
$ws.Cells.Item(82, 8).Value = 63624.75
$ws.Cells.Item(82, 9).Value = 27499.5
$ws.Cells.Item(82, 10).Value = 99750
$ws.Cells.Item(82, 11).Value = 27499.5
$ws.Cells.Item(82, 12).Value = 99750
$ws.Cells.Item(82, 13).Value = -27116.5
$ws.Cells.Item(82, 14).Value = -100516
$ws.Cells.Item(85, 8).Value = 63624.75
$ws.Cells.Item(85, 9).Value = 27499.5
$ws.Cells.Item(85, 10).Value = 99750
$ws.Cells.Item(85, 11).Value = 27499.5
$ws.Cells.Item(85, 12).Value = 99750
$ws.Cells.Item(85, 13).Value = -26173.5
$ws.Cells.Item(85, 14).Value = -102402
$ws.Cells.Item(97, 8).Value = 39999
$ws.Cells.Item(97, 9).Value = 14997
$ws.Cells.Item(97, 11).Value = 14997
$ws.Cells.Item(97, 13).Value = -14006
$ws.Cells.Item(105, 8).Value = 36517.656
$ws.Cells.Item(105, 9).Value = 49354
$ws.Cells.Item(105, 11).Value = 49354
$ws.Cells.Item(105, 13).Value = -47607
$ws.Cells.Item(140, 8).Value = 89699
$ws.Cells.Item(140, 10).Value = 89699
$ws.Cells.Item(140, 12).Value = 89699
$ws.Cells.Item(140, 14).Value = -100059

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4398.129
$ws.Cells.Item(31, 9).Value = 3319.1738
$ws.Cells.Item(31, 11).Value = 3319.1738
$ws.Cells.Item(31, 13).Value = -3024.1738
$ws.Cells.Item(34, 8).Value = 4398.129
$ws.Cells.Item(34, 9).Value = 3319.1738
$ws.Cells.Item(34, 11).Value = 3319.1738
$ws.Cells.Item(34, 13).Value = -3117.1738
$ws.Cells.Item(48, 8).Value = 66995
$ws.Cells.Item(48, 9).Value = 0
$ws.Cells.Item(48, 11).Value = 0
$ws.Cells.Item(48, 13).ClearContents()
$ws.Cells.Item(62, 8).Value = 19000
$ws.Cells.Item(62, 9).Value = 18000
$ws.Cells.Item(62, 11).Value = 18000
$ws.Cells.Item(62, 13).Value = -17376
$ws.Cells.Item(65, 8).Value = 19000
$ws.Cells.Item(65, 9).Value = 18000
$ws.Cells.Item(65, 11).Value = 90000
$ws.Cells.Item(65, 13).Value = -86880
$ws.Cells.Item(86, 8).Value = 13829.5
$ws.Cells.Item(86, 9).Value = 12139
$ws.Cells.Item(86, 10).Value = 16196.2
$ws.Cells.Item(86, 11).Value = 12139
$ws.Cells.Item(86, 12).Value = 16196.2
$ws.Cells.Item(86, 13).Value = -11016
$ws.Cells.Item(86, 14).Value = -18442.2
$ws.Cells.Item(89, 8).Value = 13829.5
$ws.Cells.Item(89, 9).Value = 12139
$ws.Cells.Item(89, 10).Value = 16196.2
$ws.Cells.Item(89, 11).Value = 60695
$ws.Cells.Item(89, 12).Value = 80981
$ws.Cells.Item(89, 13).Value = -55079
$ws.Cells.Item(89, 14).Value = -92213
$ws.Cells.Item(107, 8).Value = 100034470
$ws.Cells.Item(107, 10).Value = 2748.25
$ws.Cells.Item(107, 12).Value = 2748.25
$ws.Cells.Item(107, 14).Value = -6588.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(75, 8).Value = 1319.8
$ws.Cells.Item(75, 9).Value = 500
$ws.Cells.Item(75, 10).Value = 1524.75
$ws.Cells.Item(75, 11).Value = 1500
$ws.Cells.Item(75, 12).Value = 4574.25
$ws.Cells.Item(75, 13).Value = -502
$ws.Cells.Item(75, 14).Value = -6570.25
$ws.Cells.Item(78, 8).Value = 1319.8
$ws.Cells.Item(78, 9).Value = 500
$ws.Cells.Item(78, 10).Value = 1524.75
$ws.Cells.Item(78, 11).Value = 4500
$ws.Cells.Item(78, 12).Value = 13722.75
$ws.Cells.Item(78, 13).Value = 492
$ws.Cells.Item(78, 14).Value = -23706.75
$ws.Cells.Item(103, 8).Value = 5677.4
$ws.Cells.Item(103, 9).Value = 6966.375
$ws.Cells.Item(103, 11).Value = 20899.125
$ws.Cells.Item(103, 13).Value = -20020.125
$ws.Cells.Item(122, 8).Value = 1707.8837
$ws.Cells.Item(122, 10).Value = 2101.3667
$ws.Cells.Item(122, 12).Value = 18912.3003
$ws.Cells.Item(122, 14).Value = -23812.3003
$ws.Cells.Item(129, 8).Value = 33334614
$ws.Cells.Item(129, 9).Value = 612.7143
$ws.Cells.Item(129, 11).Value = 1838.1429
$ws.Cells.Item(129, 13).Value = 3161.8571
$ws.Cells.Item(132, 8).Value = 27279.053
$ws.Cells.Item(132, 9).Value = 684.625
$ws.Cells.Item(132, 11).Value = 6161.625
$ws.Cells.Item(132, 13).Value = -3631.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 3716.5
$ws.Cells.Item(113, 9).Value = 2766.6667
$ws.Cells.Item(113, 11).Value = 2766.6667
$ws.Cells.Item(113, 13).Value = -596.6667000000002
$ws.Cells.Item(132, 8).Value = 2486.4194
$ws.Cells.Item(132, 9).Value = 1519.6522
$ws.Cells.Item(132, 11).Value = 4558.9566
$ws.Cells.Item(132, 13).Value = -2028.9566

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 33669.707
$ws.Cells.Item(7, 9).Value = 37359
$ws.Cells.Item(7, 11).Value = 37359
$ws.Cells.Item(7, 13).Value = -37247
$ws.Cells.Item(46, 8).Value = 3144.1428
$ws.Cells.Item(46, 9).Value = 1292.8889
$ws.Cells.Item(46, 11).Value = 1292.8889
$ws.Cells.Item(46, 13).Value = -1104.8889
$ws.Cells.Item(61, 8).Value = 2299.6
$ws.Cells.Item(61, 9).Value = 2299.6
$ws.Cells.Item(61, 10).Value = 0
$ws.Cells.Item(61, 11).Value = 2299.6
$ws.Cells.Item(61, 12).Value = 0
$ws.Cells.Item(61, 13).Value = -2097.6
$ws.Cells.Item(61, 14).ClearContents()
$ws.Cells.Item(82, 8).Value = 3650.3572
$ws.Cells.Item(82, 9).Value = 4088.5
$ws.Cells.Item(82, 11).Value = 4088.5
$ws.Cells.Item(82, 13).Value = -3727.5
$ws.Cells.Item(85, 8).Value = 3650.3572
$ws.Cells.Item(85, 9).Value = 4088.5
$ws.Cells.Item(85, 11).Value = 4088.5
$ws.Cells.Item(85, 13).Value = -2840.5
$ws.Cells.Item(100, 8).Value = 3217.1924
$ws.Cells.Item(100, 9).Value = 3311.5881
$ws.Cells.Item(100, 10).Value = 3038.889
$ws.Cells.Item(100, 11).Value = 3311.5881
$ws.Cells.Item(100, 12).Value = 3038.889
$ws.Cells.Item(100, 13).Value = -2770.5881
$ws.Cells.Item(100, 14).Value = -4120.889
$ws.Cells.Item(113, 8).Value = 2299.6
$ws.Cells.Item(113, 9).Value = 2299.6
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 11).Value = 2299.6
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 13).Value = -129.5999999999999
$ws.Cells.Item(113, 14).ClearContents()
$ws.Cells.Item(126, 8).Value = 33669.707
$ws.Cells.Item(126, 9).Value = 37359
$ws.Cells.Item(126, 11).Value = 112077
$ws.Cells.Item(126, 13).Value = -109607
$ws.Cells.Item(132, 8).Value = 2622173.5
$ws.Cells.Item(132, 9).Value = 3236849.8
$ws.Cells.Item(132, 10).Value = 9800.25
$ws.Cells.Item(132, 11).Value = 9710549.399999999
$ws.Cells.Item(132, 12).Value = 29400.75
$ws.Cells.Item(132, 13).Value = -9708019.399999999
$ws.Cells.Item(132, 14).Value = -34460.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 184997.77
$ws.Cells.Item(62, 9).Value = 517577.34
$ws.Cells.Item(62, 10).Value = 3590.7273
$ws.Cells.Item(62, 11).Value = 517577.34
$ws.Cells.Item(62, 12).Value = 3590.7273
$ws.Cells.Item(62, 13).Value = -516953.34
$ws.Cells.Item(62, 14).Value = -4838.7273
$ws.Cells.Item(65, 8).Value = 184997.77
$ws.Cells.Item(65, 9).Value = 517577.34
$ws.Cells.Item(65, 10).Value = 3590.7273
$ws.Cells.Item(65, 11).Value = 2587886.7
$ws.Cells.Item(65, 12).Value = 17953.6365
$ws.Cells.Item(65, 13).Value = -2584766.7
$ws.Cells.Item(65, 14).Value = -24193.6365

